$d = $word.ActiveDocument

# The document has a "first page" header/footer that differs from the
# default (primary) header/footer - Sections(1).Headers/Footers index 1 is
# the default/primary one, index 2 is the first-page one. Each header holds
# the BTec logo picture (currently named "image1.jpg") and each footer
# holds the Pearson logo picture (currently named "image2.png"). Rename the
# pictures so the BTec logo becomes "image2.jpg" and the Pearson logo
# becomes "image1.png", in both the default and the first-page header and
# footer.

function Rename-InlinePicture($shapeRange, $newName) {
    $shapeCount = $shapeRange.InlineShapes.Count
    for ($pic = 1; $pic -le $shapeCount; $pic++) {
        $inlineShape = $shapeRange.InlineShapes($pic)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape()
    }
}

$headerCount = $d.Sections(1).Headers.Count
for ($hdrIdx = 1; $hdrIdx -le $headerCount; $hdrIdx++) {
    $header = $d.Sections(1).Headers($hdrIdx)
    if ($header.Exists) {
        Rename-InlinePicture $header.Range "image2.jpg"
    }
}

$footerCount = $d.Sections(1).Footers.Count
for ($ftrIdx = 1; $ftrIdx -le $footerCount; $ftrIdx++) {
    $footer = $d.Sections(1).Footers($ftrIdx)
    if ($footer.Exists) {
        Rename-InlinePicture $footer.Range "image1.png"
    }
}
